$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate the Danish label cells into English ---
# (Order matters for how new shared strings get appended, matching the
# original author's edit order.)
$ws.Range("B10").Value = "ID = what packet type it is"
$ws.Range("B11").Value = "Data 0-7 = the induvidual data that is being send"
$ws.Range("B21").Value = "Vbatt = Analog value 00-FF(HEX) or 0-255(Binary)"
$ws.Range("B22").Value = "At setup 00 = don't care / stay where you are"
$ws.Range("B32").Value = "Ready = is the reciver ready / Booted?"
$ws.Range("B2").Value = "Div Protocol for RC RX to TX"
$ws.Range("B23").Value = "Byte 1-4 = return current setup"

# --- Fix the alignment of the "NIU = Not In Use" note so it matches the
#     other left-aligned notes above it instead of being centered ---
$ws.Range("B12:H12").HorizontalAlignment = -4131

# --- Update the active selection ---
$ws.Range("O12").Select()
